# M4-79254 - Implement option to overrule display's Excel template in
# Export to Excel popup.
#
# The "Summary" sheet has two mirrored sections: Angle details (columns
# A:B) and Display details (columns E:F). The Angle section already has
# a "Description" row (A10/B10 -> {angledescription}); this adds the
# matching "Description" row for the Display section (E10/F10 ->
# {displaydescription}), reusing the exact look of the existing label.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# Clone the "Description" label's formatting (bold white-on-blue, wrapped)
# from A10 onto E10 so the new label matches its Angle-side counterpart.
$ws.Range("A10").Copy()
$ws.Range("E10").PasteSpecial(-4122)  # xlPasteFormats

# Populate the new Display-details "Description" row.
$ws.Range("E10").Value = "Description"
$ws.Range("F10").Value = "{displaydescription}"

# Leave the selection where the author left it when saving.
$ws.Range("A13:B13").Select() | Out-Null
